# Business Documentation update - only one KPI to do now
# Replaces the per-year "DMV_aanrijding_of_object_op_spoor" formulas (rows 2-10)
# with their hardcoded outcome values, and adds a new "Percentage" column (F)
# that expresses that count as a percentage of "Aantal_storingen".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header for the added "Percentage" KPI column
$ws.Range("F1").Value = "Percentage"

# Rows 2-10: the DMV_aanrijding_of_object_op_spoor column (C) is no longer
# computed with a formula - replace with the resulting static values.
$cValues = @{
    2  = 293
    3  = 327
    4  = 331
    5  = 290
    6  = 395
    7  = 373
    8  = 488
    9  = 521
    10 = 641
}
foreach ($row in $cValues.Keys) {
    $ws.Range("C$row").Value = $cValues[$row]
}

# Row 11-12 still compute column C via the original formula.
$ws.Range("C11:C12").Formula = "=B11/100*14.60017735"

# New "Percentage" column: C / B * 100
$ws.Range("F2").Formula = "=C2/B2*100"
$ws.Range("F2").Style = "Normal"

$ws.Range("F3:F12").Formula = "=C3/B3*100"
$ws.Range("F3:F12").Style = "Normal"

# Restore the active selection recorded in the saved workbook
$ws.Range("I7").Select()
